$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (s="2", custom date number format) from the
# last existing data row down onto the new rows so A270:A301 match the
# existing formatting exactly.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)

$data = @(
    @(270,44344,0,3,30.41979314540661),
    @(271,44345,0,3,30.41979314540661),
    @(272,44346,0,3,30.41979314540661),
    @(273,44347,0,1,10.13993104846887),
    @(274,44348,0,0,0),
    @(275,44349,0,0,0),
    @(276,44350,0,0,0),
    @(277,44351,1,1,10.13993104846887),
    @(278,44352,1,2,20.27986209693774),
    @(279,44353,0,2,20.27986209693774),
    @(280,44354,0,2,20.27986209693774),
    @(281,44355,1,3,30.41979314540661),
    @(282,44356,0,3,30.41979314540661),
    @(283,44357,0,3,30.41979314540661),
    @(284,44358,1,3,30.41979314540661),
    @(285,44359,0,2,20.27986209693774),
    @(286,44360,0,2,20.27986209693774),
    @(287,44361,0,2,20.27986209693774),
    @(288,44362,0,1,10.13993104846887),
    @(289,44363,1,2,20.27986209693774),
    @(290,44364,1,3,30.41979314540661),
    @(291,44365,0,2,20.27986209693774),
    @(292,44366,0,2,20.27986209693774),
    @(293,44367,0,2,20.27986209693774),
    @(294,44368,0,2,20.27986209693774),
    @(295,44369,0,2,20.27986209693774),
    @(296,44370,0,1,10.13993104846887),
    @(297,44371,0,0,0),
    @(298,44372,0,0,0),
    @(299,44373,1,1,10.13993104846887),
    @(300,44374,0,1,10.13993104846887),
    @(301,44375,0,1,10.13993104846887)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
